$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values are stored as text, matching the source
# data which contains thousand-separator dots and fixed decimal formatting
# that must not be reinterpreted as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.260.94'
$ws.Range("E2").Value = '  +0.14%  '
$ws.Range("D3").Value = '1.908.65'
$ws.Range("E3").Value = '  +0.12%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '307.44'
$ws.Range("E5").Value = '  -0.18%  '
$ws.Range("D7").Value = '0.5262'
$ws.Range("E7").Value = '  +1.06%  '
$ws.Range("D8").Value = '0.3811'
$ws.Range("E8").Value = '  +1.11%  '
$ws.Range("D9").Value = '0.07286'
$ws.Range("E9").Value = '  +0.13%  '
$ws.Range("D10").Value = '21.95'
$ws.Range("E10").Value = '  +3.58%  '
$ws.Range("D11").Value = '0.9019'
$ws.Range("E11").Value = '  -0.34%  '
$ws.Range("D12").Value = '0.08161'
$ws.Range("E12").Value = '  -3.41%  '
$ws.Range("D13").Value = '95.98'
$ws.Range("E13").Value = '  -1.02%  '
$ws.Range("D14").Value = '5.360'
$ws.Range("E14").Value = '  +1.14%  '
$ws.Range("D15").Value = '1.451.97'
$ws.Range("E15").Value = '  -23.87%  '
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  -0.09%  '
$ws.Range("D17").Value = '0.000008656'
$ws.Range("E17").Value = '  -0.22%  '
$ws.Range("E18").Value = '  +1.45%  '
$ws.Range("E19").Value = '  -0.05%  '
$ws.Range("D20").Value = '27.291.96'
$ws.Range("E20").Value = '  +0.13%  '
$ws.Range("D21").Value = '5.117'
$ws.Range("E21").Value = '  +0.38%  '
$ws.Range("E22").Value = '  +1.60%  '
$ws.Range("E23").Value = '  +1.05%  '
$ws.Range("D24").Value = '149.81'
$ws.Range("E24").Value = '  +2.01%  '
$ws.Range("D25").Value = '2.307'
$ws.Range("E25").Value = '  -0.99%  '
$ws.Range("D26").Value = '18.26'
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("D27").Value = '1.736'
$ws.Range("E27").Value = '  -1.15%  '
$ws.Range("D28").Value = '116.94'
$ws.Range("E28").Value = '  +1.50%  '
$ws.Range("D29").Value = '4.847'
$ws.Range("E29").Value = '  +0.40%  '
$ws.Range("D30").Value = '4.849'
$ws.Range("E30").Value = '  -1.30%  '
$ws.Range("D31").Value = '0.09250'
$ws.Range("E31").Value = '  -0.37%  '
$ws.Range("D32").Value = '0.8301'
$ws.Range("E32").Value = '  +4.26%  '
$ws.Range("D33").Value = '0.05067'
$ws.Range("E33").Value = '  -0.25%  '
$ws.Range("D34").Value = '1.228'
$ws.Range("E34").Value = '  -1.33%  '
$ws.Range("D35").Value = '2.993'
$ws.Range("E35").Value = '  +1.45%  '
$ws.Range("D36").Value = '3.351'
$ws.Range("E36").Value = '  -2.28%  '
$ws.Range("D37").Value = '2.717'
$ws.Range("E37").Value = '  +4.95%  '
$ws.Range("D38").Value = '0.5796'
$ws.Range("E38").Value = '  -0.05%  '
$ws.Range("D39").Value = '0.02002'
$ws.Range("E39").Value = '  -0.31%  '
$ws.Range("D40").Value = '1.076'
$ws.Range("E40").Value = '  -0.02%  '
$ws.Range("D41").Value = '9.230'
$ws.Range("E41").Value = '  +1.73%  '
$ws.Range("D42").Value = '6.594'
$ws.Range("E42").Value = '  -0.38%  '
$ws.Range("E43").Value = '  -0.13%  '
$ws.Range("E44").Value = '  +0.15%  '
$ws.Range("D45").Value = '0.4927'
$ws.Range("E45").Value = '  +0.97%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '10.19'
$ws.Range("E46").Value = '  +0.28%  '
$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").Value = '1.000'
$ws.Range("E47").Value = '  -0.14%  '
$ws.Range("E48").Value = '  +0.53%  '
$ws.Range("D49").Value = '39.03'
$ws.Range("E49").Value = '  +3.37%  '
$ws.Range("D50").Value = '0.06172'
$ws.Range("E50").Value = '  +3.52%  '
$ws.Range("E51").Value = '  +0.64%  '
